$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.482117666666667
$ws.Range("N2").Value = 4.446353
$ws.Range("O2").Value = 0.1401829251394648
$ws.Range("P2").Value = 0.1401829251394648
$ws.Range("Q2").Value = 4.285550149715777
$ws.Range("R2").Value = 38.569951347442
$ws.Range("S2").Value = 0.01701246600963427
$ws.Range("T2").Value = 0.01701246600963427

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.365790333333333
$ws.Range("N3").Value = 7.097371
$ws.Range("O3").Value = 0.2237632116883227
$ws.Range("P3").Value = 0.2237632116883226
$ws.Range("Q3").Value = 6.840693789188222
$ws.Range("R3").Value = 61.566244102694
$ws.Range("S3").Value = 0.02715568981933373
$ws.Range("T3").Value = 0.02715568981933372

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.724832333333333
$ws.Range("N4").Value = 20.174497
$ws.Range("O4").Value = 0.6360538631722126
$ws.Range("P4").Value = 0.6360538631722126
$ws.Range("Q4").Value = 19.44488407438422
$ws.Range("R4").Value = 175.003956669458
$ws.Range("S4").Value = 0.07719088980878677
$ws.Range("T4").Value = 0.07719088980878677

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.482117666666667
$ws.Range("N5").Value = 4.446353
$ws.Range("O5").Value = 0.1401829251394648
$ws.Range("P5").Value = 0.1401829251394648
$ws.Range("Q5").Value = 17.84533352322411
$ws.Range("R5").Value = 160.608001709017
$ws.Range("S5").Value = 0.07084111009984841
$ws.Range("T5").Value = 0.07084111009984839

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.365790333333333
$ws.Range("N6").Value = 7.097371
$ws.Range("O6").Value = 0.2237632116883227
$ws.Range("P6").Value = 0.2237632116883226
$ws.Range("Q6").Value = 28.48513211457989
$ws.Range("R6").Value = 256.3661890312191
$ws.Range("S6").Value = 0.1130782104863179
$ws.Range("T6").Value = 0.1130782104863179

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.724832333333333
$ws.Range("N7").Value = 20.174497
$ws.Range("O7").Value = 0.6360538631722126
$ws.Range("P7").Value = 0.6360538631722126
$ws.Range("Q7").Value = 80.9698707296259
$ws.Range("R7").Value = 728.728836566633
$ws.Range("S7").Value = 0.3214283173616807
$ws.Range("T7").Value = 0.3214283173616807

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.894099000000001
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.482117666666667
$ws.Range("N8").Value = 4.446353
$ws.Range("O8").Value = 0.1401829251394648
$ws.Range("P8").Value = 0.1401829251394648
$ws.Range("Q8").Value = 13.18210125698233
$ws.Range("R8").Value = 118.638911312841
$ws.Range("S8").Value = 0.05232934902998215
$ws.Range("T8").Value = 0.05232934902998213

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.894099000000001
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.365790333333333
$ws.Range("N9").Value = 7.097371
$ws.Range("O9").Value = 0.2237632116883227
$ws.Range("P9").Value = 0.2237632116883226
$ws.Range("Q9").Value = 21.04157343790967
$ws.Range("R9").Value = 189.374160941187
$ws.Range("S9").Value = 0.08352931138267101
$ws.Range("T9").Value = 0.08352931138267099

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.894099000000001
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.724832333333333
$ws.Range("N10").Value = 20.174497
$ws.Range("O10").Value = 0.6360538631722126
$ws.Range("P10").Value = 0.6360538631722126
$ws.Range("Q10").Value = 59.81132453106766
$ws.Range("R10").Value = 538.301920779609
$ws.Range("S10").Value = 0.2374346560017451
$ws.Range("T10").Value = 0.2374346560017451
